$d = $word.ActiveDocument

$olds = @(
    "<x0>",
    "<g0>oremlay ipsumhay olorday itsay amethay, onsectetuercay adipiscinghay elithay. aecenasmay orttitorpay onguecay assamay. uscefay osuerepay, agnamay edsay ulvinarpay ultricieshay, uruspay ectuslay alesuadamay iberolay, itsay amethay ommodocay agnamay eroshay uisqay urnahay.</g1>",
    "<g0>uncnay iverravay imperdiethay enimhay. uscefay esthay. ivamusvay ahay ellustay.</g1>",
    "<g0>ellentesquepay </g1><g2>abitanthay orbimay istiquetray</g3><g4> enectussay ethay etusnay ethay alesuadamay amesfay achay urpistay egestashay. oinpray aretraphay onummynay edepay. </g5><x6><x7><x8><x9><x10><x11><x12><x13><g14>aurishay ethay orcihay.</g15>",
    "<g0>Aeneanhay ecnay oremlay. </g1><g2>Inhay orttitorpay. onecday aoreetlay onummynay auguehay.</g3>",
    "<g0>uspendissesay uiday uruspay, elerisquescay athay, </g1><g2>ulputatevay</g3><g4> itaevay, </g5><g6>etiumpray</g7><g8> attismay, uncnay. </g9><g10>aurismay egethay equenay athay emsay enenatisvay eleifendhay. Uthay onummynay.</g11>",
    "<g0>uscefay aliquethay edepay onnay edepay. </g1><g2>uspendissesay apibusday oremlay ellentesquepay agnamay. Integerhay ullanay.</g3>",
    "<g0>onecday anditblay eugiatfay igulalay. </g1><g2>onecday endrerithay</g3><g4>, elisfay ethay imperdiethay euismodhay, uruspay ipsumhay etiumpray etusmay, inhay acinialay ullanay islnay egethay apiensay. onecday uthay esthay inhay ectuslay </g5><g6>onsequatcay</g7><g8> onsequatcay.</g9>",
    "<g0>Etiamhay egethay uiday. </g1><g2>Aliquamhay erathay olutpatvay.</g3><g4> edsay athay oremlay inhay uncnay ortapay istiquetray.</g5>",
    "<g0>oinpray ecnay auguehay. </g1><g2>uisqueqay</g3><g4> aliquamhay </g5><g6>emportay</g7><g8> agnamay. ellentesquepay abitanthay orbimay istiquetray enectussay ethay etusnay ethay alesuadamay amesfay achay urpistay egestashay.</g9>",
    "<g0>uncnay achay agnamay. aecenasmay odiohay olorday, </g1><g2>ulputatevay </g3><g4>elvay, auctorhay achay, accumsanhay idhay, elisfay. </g5><g6>eehahyay</g7>ellentesquepay<x8><g9> ursuscay agittissay elisfay.</g10>"
)

$news = @(
    "<g0><x1></g2><x3>",
    "<g0><g1><x2></g3></g4><g5>oremlay ipsumhay olorday itsay amethay, onsectetuercay adipiscinghay elithay. aecenasmay orttitorpay onguecay assamay. uscefay osuerepay, agnamay edsay ulvinarpay ultricieshay, uruspay ectuslay alesuadamay iberolay, itsay amethay ommodocay agnamay eroshay uisqay urnahay.</g6>",
    "<g0><g1><x2></g3></g4><g5>uncnay iverravay imperdiethay enimhay. uscefay esthay. ivamusvay ahay ellustay.</g6>",
    "<g0><g1><x2></g3></g4><g5>ellentesquepay </g6><g7>abitanthay orbimay istiquetray</g8><g9> enectussay ethay etusnay ethay alesuadamay amesfay achay urpistay egestashay. oinpray aretraphay onummynay edepay. </g10><x11><x12><x13><x14><x15><x16><x17><x18><g19>aurishay ethay orcihay.</g20>",
    "<g0><g1><x2><x3></g4></g5><g6>Aeneanhay ecnay oremlay. </g7><g8>Inhay orttitorpay. onecday aoreetlay onummynay auguehay.</g9>",
    "<g0><g1><x2><x3></g4></g5><g6>uspendissesay uiday uruspay, elerisquescay athay, </g7><g8>ulputatevay</g9><g10> itaevay, </g11><g12>etiumpray</g13><g14> attismay, uncnay. </g15><g16>aurismay egethay equenay athay emsay enenatisvay eleifendhay. Uthay onummynay.</g17>",
    "<g0><g1><x2></g3></g4><g5>uscefay aliquethay edepay onnay edepay. </g6><g7>uspendissesay apibusday oremlay ellentesquepay agnamay. Integerhay ullanay.</g8>",
    "<g0><g1><x2></g3></g4><g5>onecday anditblay eugiatfay igulalay. </g6><g7>onecday endrerithay</g8><g9>, elisfay ethay imperdiethay euismodhay, uruspay ipsumhay etiumpray etusmay, inhay acinialay ullanay islnay egethay apiensay. onecday uthay esthay inhay ectuslay </g10><g11>onsequatcay</g12><g13> onsequatcay.</g14>",
    "<g0><g1><x2></g3></g4><g5>Etiamhay egethay uiday. </g6><g7>Aliquamhay erathay olutpatvay.</g8><g9> edsay athay oremlay inhay uncnay ortapay istiquetray.</g10>",
    "<g0><g1><x2></g3></g4><g5>oinpray ecnay auguehay. </g6><g7>uisqueqay</g8><g9> aliquamhay </g10><g11>emportay</g12><g13> agnamay. ellentesquepay abitanthay orbimay istiquetray enectussay ethay etusnay ethay alesuadamay amesfay achay urpistay egestashay.</g14>",
    "<g0><g1><x2></g3></g4><g5>uncnay achay agnamay. aecenasmay odiohay olorday, </g6><g7>ulputatevay </g8><g9>elvay, auctorhay achay, accumsanhay idhay, elisfay. </g10><g11>eehahyay</g12>ellentesquepay<x13><g14> ursuscay agittissay elisfay.</g15>"
)

for ($i = 0; $i -lt $olds.Length; $i++) {
    $para = $d.Paragraphs.Item($i + 1)
    $rng = $para.Range
    $ok = $rng.Find.Execute($olds[$i], $true, $false, $false, $false, $false, $true, 1, $false, $news[$i], 2)
    if (-not $ok) {
        Write-Output "WARNING: replacement $($i + 1) did not find its target text"
    }
}

Write-Output "Done: retagged $($olds.Length) paragraphs"
